$d = $word.ActiveDocument

$replacements = @(
    @("80×63=", "93×95="),
    @("78×24=", "45×95="),
    @("14×23=", "13×99="),
    @("13×67=", "50×24="),
    @("50×58=", "26×69="),
    @("46×99=", "46×27="),
    @("29×31=", "65×55="),
    @("66×59=", "75×41="),
    @("11×34=", "63×24="),
    @("20×78=", "64×69="),
    @("71×17=", "52×18="),
    @("34×57=", "76×83="),
    @("79×98=", "71×47="),
    @("63×45=", "72×73="),
    @("94×20=", "65×23="),
    @("57×73=", "37×20="),
    @("50×48=", "57×14="),
    @("19×57=", "66×88="),
    @("29×44=", "38×48="),
    @("78×69=", "32×73="),
    @("68×40=", "85×73="),
    @("55×69=", "36×79="),
    @("14×77=", "34×11="),
    @("11×38=", "97×40="),
    @("92×28=", "40×90=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
